# Auto-generated edit script: updates cryptos list (prices / 1h volume %, and
# a couple of row re-orderings for Chainlink/WrappedEther and ThetaToken/VeChain)
# to match the "Thu Mar 28 14:25:58 UTC 2024" GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $NewValue)
    $range = $ws.Range($Cell)
    $origStyle = $range.Style
    if ($NewValue -match '^\s*[-+]?(\d+\.?\d*|\.\d+)\s*$') {
        # Guard against Excel auto-converting a numeric-looking string into a
        # real number (the source data keeps these as plain text cells).
        $range.Value = "'" + $NewValue
    } else {
        $range.Value = $NewValue
    }
    # Re-apply the original style so a forced text entry ("quote prefix")
    # doesn't leave the cell's formatting different from before.
    $range.Style = $origStyle
}

Set-TextValue 'D2' '70.942.26'
Set-TextValue 'E2' '  +2.06%  '
Set-TextValue 'D3' '3.579.35'
Set-TextValue 'E3' '  +1.37%  '
Set-TextValue 'E4' '  +0.16%  '
Set-TextValue 'D5' '585.72'
Set-TextValue 'E5' '  +2.68%  '
Set-TextValue 'D6' '185.81'
Set-TextValue 'E6' '  +0.96%  '
Set-TextValue 'D7' '3.571.37'
Set-TextValue 'D8' '0.621'
Set-TextValue 'E8' '  +1.05%  '
Set-TextValue 'E9' '  +0.05%  '
Set-TextValue 'D10' '0.219'
Set-TextValue 'E10' '  +19.65%  '
Set-TextValue 'D11' '0.651'
Set-TextValue 'E11' '  +0.77%  '
Set-TextValue 'D12' '54.17'
Set-TextValue 'E12' '  +0.47%  '
Set-TextValue 'D13' '0.0000322'
Set-TextValue 'E13' '  +7.71%  '
Set-TextValue 'D14' '9.52'
Set-TextValue 'E14' '  +0.42%  '
Set-TextValue 'D15' '4.146.64'
Set-TextValue 'E15' '  +1.24%  '
Set-TextValue 'D16' '70.806.84'
Set-TextValue 'E16' '  +2.36%  '
Set-TextValue 'B17' 'WrappedEther'
Set-TextValue 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '3.589.70'
Set-TextValue 'E17' '  +2.17%  '
Set-TextValue 'B18' 'Chainlink'
Set-TextValue 'C18' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D18' '19.22'
Set-TextValue 'E18' '  -0.61%  '
Set-TextValue 'D19' '573.67'
Set-TextValue 'E19' '  +13.38%  '
Set-TextValue 'D20' '12.35'
Set-TextValue 'E20' '  +0.14%  '
Set-TextValue 'D21' '0.120'
Set-TextValue 'E21' '  +0.63%  '
Set-TextValue 'E22' '  -2.19%  '
Set-TextValue 'D23' '17.44'
Set-TextValue 'E23' '  -10.76%  '
Set-TextValue 'D24' '5.08'
Set-TextValue 'E24' '  +3.43%  '
Set-TextValue 'D25' '4.60'
Set-TextValue 'E25' '  +6.37%  '
Set-TextValue 'D26' '94.78'
Set-TextValue 'E26' '  +0.79%  '
Set-TextValue 'D27' '11.29'
Set-TextValue 'E27' '  -0.31%  '
Set-TextValue 'D28' '2.92'
Set-TextValue 'E28' '  -0.31%  '
Set-TextValue 'D29' '9.09'
Set-TextValue 'E29' '  -1.01%  '
Set-TextValue 'D30' '32.28'
Set-TextValue 'E30' '  +2.69%  '
Set-TextValue 'D31' '7.20'
Set-TextValue 'E31' '  -4.69%  '
Set-TextValue 'D32' '12.25'
Set-TextValue 'E32' '  -1.35%  '
Set-TextValue 'D33' '0.114'
Set-TextValue 'E33' '  -0.44%  '
Set-TextValue 'D34' '64.13'
Set-TextValue 'E34' '  -1.64%  '
Set-TextValue 'E35' '  +5.83%  '
Set-TextValue 'D36' '551.62'
Set-TextValue 'E36' '  -3.28%  '
Set-TextValue 'D37' '0.413'
Set-TextValue 'E37' '  +2.83%  '
Set-TextValue 'D38' '0.0₃0812'
Set-TextValue 'E38' '  +4.53%  '
Set-TextValue 'D39' '37.49'
Set-TextValue 'E39' '  -1.27%  '
Set-TextValue 'E40' '  -0.04%  '
Set-TextValue 'D41' '3.500.09'
Set-TextValue 'E41' '  +10.37%  '
Set-TextValue 'D42' '3.19'
Set-TextValue 'E42' '  -0.43%  '
Set-TextValue 'D43' '3.45'
Set-TextValue 'E43' '  +1.77%  '
Set-TextValue 'E44' '  +1.54%  '
Set-TextValue 'D45' '3.51'
Set-TextValue 'E45' '  -0.62%  '
Set-TextValue 'B46' 'VeChain'
Set-TextValue 'C46' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D46' '0.0445'
Set-TextValue 'E46' '  +0.07%  '
Set-TextValue 'B47' 'ThetaToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D47' '2.96'
Set-TextValue 'E47' '  -0.83%  '
Set-TextValue 'D48' '9.38'
Set-TextValue 'E48' '  +0.93%  '
Set-TextValue 'D49' '0.137'
Set-TextValue 'E49' '  +2.05%  '
Set-TextValue 'D50' '0.998'
Set-TextValue 'E50' '  +0.12%  '
Set-TextValue 'D51' '1.45'
Set-TextValue 'E51' '  -0.65%  '
